$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.160.17"
$ws.Range("E2").Value = "  +1.63%  "

$ws.Range("D3").Value = "2.579.30"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.63"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("E6").Value = "  -2.50%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -7.77%  "

$ws.Range("D9").Value = "2.583.15"
$ws.Range("E9").Value = "  +0.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.61"

$ws.Range("E11").Value = "  +0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").Value = "3.034.07"
$ws.Range("E14").Value = "  +0.30%  "

$ws.Range("D15").Value = "60.205.15"
$ws.Range("E15").Value = "  +1.67%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.43"
$ws.Range("E16").Value = "  -1.56%  "

$ws.Range("E17").Value = "  +1.92%  "

$ws.Range("D18").Value = "2.585.26"
$ws.Range("E18").Value = "  +0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.05"
$ws.Range("E20").Value = "  +2.97%  "

$ws.Range("E21").Value = "  +0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  +1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.04"
$ws.Range("E24").Value = "  +0.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  +1.60%  "

$ws.Range("E26").Value = "  -0.57%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("E27").Value = "  -0.65%  "

$ws.Range("E28").Value = "  +2.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.32"
$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.41"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("E33").Value = "  -0.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.74"
$ws.Range("E34").Value = "  +4.69%  "

$ws.Range("E35").Value = "  +1.58%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.856"
$ws.Range("E37").Value = "  +16.97%  "

$ws.Range("E38").Value = "  -0.11%  "

$ws.Range("E39").Value = "  +0.12%  "

$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.83"
$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.23"
$ws.Range("E42").Value = "  +1.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0997"
$ws.Range("E43").Value = "  -1.96%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0558"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("E45").Value = "  -1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("E46").Value = "  -0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.73"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.88"
$ws.Range("E48").Value = "  +2.54%  "

$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("E50").Value = "  +0.16%  "

$ws.Range("D51").Value = "1.999.21"
$ws.Range("E51").Value = "  +0.89%  "
